$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.695620666666667
$ws.Range("N2").Value = 23.086862
$ws.Range("O2").Value = 0.4976976897997125
$ws.Range("P2").Value = 0.4976976897997126
$ws.Range("Q2").Value = 0.5234304656777778
$ws.Range("R2").Value = 4.7108741911
$ws.Range("S2").Value = 0.4976976897997125
$ws.Range("T2").Value = 0.4976976897997126

# Row 3 (Target cluster: FAPs)
$ws.Range("O3").Value = 0.228990810419744
$ws.Range("P3").Value = 0.228990810419744
$ws.Range("S3").Value = 0.228990810419744
$ws.Range("T3").Value = 0.228990810419744

# Row 4 (Target cluster: MuSCs)
$ws.Range("M4").Value = 3.610968333333334
$ws.Range("N4").Value = 10.832905
$ws.Range("O4").Value = 0.2335315987213747
$ws.Range("P4").Value = 0.2335315987213747
$ws.Range("Q4").Value = 0.2456060294722222
$ws.Range("R4").Value = 2.21045426525
$ws.Range("S4").Value = 0.2335315987213747
$ws.Range("T4").Value = 0.2335315987213747

# Row 5 (Target cluster: Resolving-Mac)
$ws.Range("M5").Value = 0.6150943333333334
$ws.Range("N5").Value = 1.845283
$ws.Range("O5").Value = 0.03977990105916876
$ws.Range("P5").Value = 0.03977990105916877
$ws.Range("Q5").Value = 0.0418366662388889
$ws.Range("R5").Value = 0.37652999615
$ws.Range("S5").Value = 0.03977990105916876
$ws.Range("T5").Value = 0.03977990105916877
